# Success Analysis.xlsx update:
#  - Add Table1 (Sheet1) a new trailing column "Column1"
#  - Add a new worksheet "Full Run" (after Sheet3) with Run/Time/Points data
#  - Update selections / active sheet to match the new state

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "Full Run" worksheet at the end of the workbook
#    (done first so the new shared strings "Time"/"Points" are allocated
#    before Table1's new "Column1" header, matching the authored order)
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$fullRun = $wb.Worksheets.Add($null, $lastSheet)
$fullRun.Name = "Full Run"

$fullRun.Range("A1").Value = "Run"
$fullRun.Range("B1").Value = "Time"
$fullRun.Range("C1").Value = "Points"

$fullRun.Range("A2").Value = 1
$fullRun.Range("B2").Value = 0.10833333333333334
$fullRun.Range("B2").NumberFormat = "h:mm"
$fullRun.Range("C2").Value = 260

$fullRun.Range("A3").Value = 2
$fullRun.Range("B3").Value = 0.1125
$fullRun.Range("B3").NumberFormat = "h:mm"
$fullRun.Range("C3").Value = 355

$fullRun.Range("A4").Value = 3
$fullRun.Range("B4").Value = 0.10416666666666667
$fullRun.Range("B4").NumberFormat = "h:mm"
$fullRun.Range("C4").Value = 285

$fullRun.Range("A5").Value = 4
$fullRun.Range("B5").Value = 0.095138888888888884
$fullRun.Range("B5").NumberFormat = "h:mm"
$fullRun.Range("C5").Value = 295

# ---------------------------------------------------------------------------
# 2. Sheet1 / Table1: append a new "Column1" column to the table
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$table1 = $ws1.ListObjects.Item("Table1")
[void]$table1.ListColumns.Add()
$ws1.Range("H1").Value = "Column1"
$ws1.Columns.Item(8).ColumnWidth = 11.25

[void]$ws1.Range("V21").Select()

# ---------------------------------------------------------------------------
# 3. Sheet3: keep its data, just move the selection (it stops being the
#    active tab once "Full Run" is added)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
[void]$ws3.Range("Q28").Select()

# ---------------------------------------------------------------------------
# 4. Final selection + activation on the new sheet
# ---------------------------------------------------------------------------
[void]$fullRun.Range("F12").Select()
$fullRun.Activate()
